# "added a few more" - append 10 more batal/meaning/source rows to the table,
# plus a few trailing blank (but styled) placeholder rows, matching the
# established row layout/styling of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the existing "empty placeholder" styling (s="2", center aligned)
# down through row 74 before filling in the new data, so the freshly
# created rows pick up the same formatting as the rest of the table.
$ws.Range("A62:C64").Copy()
$ws.Range("A65:C74").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# New batal / meaning pairs (column C repeats the constant "source" value
# used throughout the sheet). Use Value2 here: reading back .Value from a
# cell in this environment doesn't resolve to the underlying scalar.
$source = $ws.Cells.Item(2, 3).Value2

$ws.Cells.Item(62, 1).Value = "بور وکمیت پہ کاش ءِ چرگ ءَ نہ ٹہیت۔"
$ws.Cells.Item(62, 2).Value = "کاش ہما بُرّناکیں کاہ ئے، کاش ءَ پہ کُڈّک ءُ پَل ءَ کارمز کناں، نیکہ پہ اولاک کدیم ءَ"
$ws.Cells.Item(62, 3).Value = $source

$ws.Cells.Item(63, 1).Value = "بور ءِ نالاں بہ خچر ءَ بِر مہ جن۔"
$ws.Cells.Item(63, 2).Value = "بور ءِ نال پہ بور ءَ اَنت، خچر ءَ بِراش نہ کناں"
$ws.Cells.Item(63, 3).Value = $source

$ws.Cells.Item(64, 1).Value = "بہ جنگ ءَ بے دل مہ بو۔"
$ws.Cells.Item(64, 2).Value = "جنگ ءِ تہا بزدلی تاوان دنت،‌ وتارا بے‌‌ دل‌ مہ کن"
$ws.Cells.Item(64, 3).Value = $source

$ws.Cells.Item(65, 1).Value = "بہ چپ ءُ چوٹیں براس ءَ دل مہ بند۔"
$ws.Cells.Item(65, 2).Value = "آ براس کہ دھو کہ باز اِنت، لالچی اِنت آئی پہ اوست مہ بند"
$ws.Cells.Item(65, 3).Value = $source

$ws.Cells.Item(66, 1).Value = "بہ‌رہ ءِ سرا ءُ کَور ءِ گُور ءَ نپاد مہ کن۔"
$ws.Cells.Item(66, 2).Value = "راہ ءُ کور ہنچیں جاہ اَنت کہ اود ءَ اڈّکنگی نہ اِنت"
$ws.Cells.Item(66, 3).Value = $source

# Note: for this particular row the "meaning" text (column B) was entered
# before the "batal" text (column A), so it claims the earlier shared-string
# slot.
$ws.Cells.Item(67, 2).Value = "آدیوان ءَ پترلگّیت جیڑ ہانی تو جیل ءَ گڑا‌شور ءُ سلاہ بیت اُود ءَ بائدیں کوکار ءُ جاک مہ بیت، توجیل کنگ بہ بیت"
$ws.Cells.Item(67, 1).Value = "بِہ سُہل ءُ سلہ ءِ نیام ءَ  غوغا مہ کن۔"
$ws.Cells.Item(67, 3).Value = $source

$ws.Cells.Item(68, 1).Value = "بِہ صید ءُ شکار ءِ نیام ءَ ہاہا مہ کن۔"
$ws.Cells.Item(68, 2).Value = "شکار ءَ کہ ردئے گڑا کوکار مہ کن، شکار پہ ہاموشی ءَ بیت"
$ws.Cells.Item(68, 3).Value = $source

$ws.Cells.Item(69, 1).Value = "بہ ہرکار ءَ وبہ ہرمار ءَ ہلگر مہ کن۔"
$ws.Cells.Item(69, 2).Value = "بزاں ہر کارے ءِ تہاوت جہد بہ کن چمدارمہ بو"
$ws.Cells.Item(69, 3).Value = $source

$ws.Cells.Item(70, 1).Value = "بہارگہ پہ بہار، بہار پہ ڈگار۔"
$ws.Cells.Item(70, 2).Value = "آپ ءُ ہئور بیت گڑا بہارگہ بیت، ڈگار پہ کِشت ءُ کشار"
$ws.Cells.Item(70, 3).Value = $source

$ws.Cells.Item(71, 1).Value = "بیکارو فضول نہ خداوش انت نہ رسوُل ؐ۔"
$ws.Cells.Item(71, 2).Value = "بے کار ءُ ناشریں مردم ءَ راکس دوست نہ داریت"
$ws.Cells.Item(71, 3).Value = $source

# Leave rows 72:74 as blank (but styled) placeholder rows, matching the
# original tail of the sheet.

$ws.Cells.Item(71, 2).Select()
